$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.096.01'
$ws.Range('E2').Value = '  -1.12%  '

$ws.Range('D3').Value = '3.102.99'
$ws.Range('E3').Value = '  -1.72%  '

$ws.Range('E4').Value = '  -0.54%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.70'
$ws.Range('E5').Value = '  +0.95%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.00'
$ws.Range('E6').Value = '  +2.33%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.27%  '

$ws.Range('E8').Value = '  +0.38%  '

$ws.Range('D9').Value = '3.103.16'
$ws.Range('E9').Value = '  -1.64%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').Value = '  -2.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.96'
$ws.Range('E11').Value = '  -0.61%  '

$ws.Range('E12').Value = '  -2.93%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000240'
$ws.Range('E13').Value = '  -3.68%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.01'
$ws.Range('E14').Value = '  -4.51%  '

$ws.Range('E15').Value = '  -0.83%  '

$ws.Range('D16').Value = '3.610.28'
$ws.Range('E16').Value = '  -1.92%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.22'
$ws.Range('E17').Value = '  -1.49%  '

$ws.Range('D18').Value = '63.965.16'
$ws.Range('E18').Value = '  -0.79%  '

$ws.Range('D19').Value = '3.101.23'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '481.63'
$ws.Range('E20').Value = '  +1.30%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.51'
$ws.Range('E21').Value = '  -3.49%  '

$ws.Range('E22').Value = '  -4.94%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.58'
$ws.Range('E23').Value = '  -1.53%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.49'
$ws.Range('E24').Value = '  +3.10%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.61'
$ws.Range('E25').Value = '  -0.98%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.93'
$ws.Range('E26').Value = '  -4.83%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.81'
$ws.Range('E27').Value = '  +7.84%  '

$ws.Range('E28').Value = '  -0.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.66'
$ws.Range('E29').Value = '  +3.04%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.70'
$ws.Range('E30').Value = '  -1.49%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.21'
$ws.Range('E31').Value = '  -1.23%  '

$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.60%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.27'
$ws.Range('E34').Value = '  -2.01%  '

$ws.Range('D35').Value = '0.0₃0846'
$ws.Range('E35').Value = '  -4.08%  '

$ws.Range('E36').Value = '  +0.91%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.05'
$ws.Range('E37').Value = '  -2.99%  '

$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.31'
$ws.Range('E38').Value = '  -6.98%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.25'
$ws.Range('E39').Value = '  -3.31%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '51.04'
$ws.Range('E40').Value = '  -0.71%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.24'
$ws.Range('E41').Value = '  -1.47%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '446.21'
$ws.Range('E42').Value = '  -4.82%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.292'
$ws.Range('E43').Value = '  -3.48%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0366'
$ws.Range('E44').Value = '  -4.34%  '

$ws.Range('E45').Value = '  +0.90%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.25'
$ws.Range('E46').Value = '  +3.52%  '

$ws.Range('D47').Value = '2.842.75'
$ws.Range('E47').Value = '  -2.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.27'
$ws.Range('E48').Value = '  +0.69%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '26.12'
$ws.Range('E49').Value = '  +0.37%  '

$ws.Range('E50').Value = '  +0.01%  '

$ws.Range('E51').Value = '  -2.29%  '

